$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header strings: "<name>_old" -> "<name>_FV2404" and
# "<name>_new" -> "<name>_FV2410" (the AHB file pair changed).
$null = $ws.Range("A1:J1").Replace("_old", "_FV2404")
$null = $ws.Range("L1:U1").Replace("_new", "_FV2410")

# Turn the data range into a real Excel table ("Table1") covering A1:U70.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
